{"js": "// The document has two occurrences of the literal text \"PF\":\n//   1. \"(Examen y PF)\"        -> bold, orange (accent2 / #ED7D31) run   <- the one we must change\n//   2. \"(Correcci\u00f3n de PF)\"   -> plain, black run                       <- must stay untouched\n//\n// The edit changes \"PF\" to \"Cuest.\" only in the first occurrence, turning\n// \"(Examen y PF)\" into \"(Examen y Cuest.)\". The Word \"last edit\" bookmark\n// (_GoBack) is relocated by Word to sit right after the newly typed text,\n// so we remove the old bookmark (it previously sat elsewhere, between\n// \", Fec.\" and \" Lim.\") and re-insert it immediately after \"Cuest.\".\n\nconst body = context.document.body;\n\n// Locate every run containing \"PF\" and inspect formatting to find the\n// bold / orange (#ED7D31) occurrence inside \"(Examen y PF)\".\nconst results = body.search(\"PF\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text,font/color,font/bold\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < results.items.length; i++) {\n  const item = results.items[i];\n  const color = (item.font.color || \"\").toUpperCase();\n  if (item.font.bold === true && color === \"#ED7D31\") {\n    target = item;\n    break;\n  }\n}\n\n// Fall back to \"the only bold match\" if the color comparison above (e.g.\n// due to a differently-formatted color string) didn't pin down a result.\nif (!target) {\n  const boldMatches = results.items.filter((item) => item.font.bold === true);\n  if (boldMatches.length === 1) {\n    target = boldMatches[0];\n  }\n}\n\nif (!target) {\n  throw new Error('Could not find the target \"PF\" run (bold, #ED7D31).');\n}\n\n// Replace \"PF\" with \"Cuest.\" in place; the run's existing formatting\n// (bold + accent2 orange color) carries over automatically.\ntarget.insertText(\"Cuest.\", Word.InsertLocation.replace);\nawait context.sync();\n\n// Move the \"_GoBack\" bookmark so it ends up right after the text we just\n// typed (matching Word's own behavior of tracking the last edit point).\nconst oldGoBack = context.document.getBookmarkRangeOrNullObject(\"_GoBack\");\noldGoBack.load(\"isNullObject\");\nawait context.sync();\n\nif (!oldGoBack.isNullObject) {\n  context.document.deleteBookmark(\"_GoBack\");\n}\n\nconst afterEdit = target.getRange(Word.RangeLocation.after);\nafterEdit.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# The document has two occurrences of the literal text \"PF\":\n#   1. \"(Examen y PF)\"        -> bold, orange (accent2 / #ED7D31) run   <- the one we must change\n#   2. \"(Correcci\u00f3n de PF)\"   -> plain, non-bold, black run              <- must stay untouched\n#\n# The edit changes \"PF\" to \"Cuest.\" only in the first occurrence, turning\n# \"(Examen y PF)\" into \"(Examen y Cuest.)\". Word also relocates the\n# \"_GoBack\" bookmark (which marks the last edit location) so that it sits\n# right after the newly typed text; we replicate that by removing the old\n# bookmark (previously between \", Fec.\" and \" Lim.\") and re-adding it\n# immediately after \"Cuest.\".\n\n$d = $word.ActiveDocument\n\n$found = $false\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"PF\"\n$find.Font.Bold = $true\n$find.MatchCase = $true\n$find.Forward = $true\n$find.Wrap = 0\n\nif ($find.Execute()) {\n    $found = $true\n    $rng = $find.Parent\n} else {\n    # Fallback: scan every \"PF\" occurrence and pick the bold one (in case\n    # the formatted Find above doesn't pin it down for some reason).\n    $find2 = $d.Content.Find\n    $find2.ClearFormatting()\n    $find2.Text = \"PF\"\n    $find2.MatchCase = $true\n    $find2.Forward = $true\n    $find2.Wrap = 0\n    while ($find2.Execute()) {\n        if ($find2.Parent.Font.Bold -eq -1) {\n            $found = $true\n            $rng = $find2.Parent\n            break\n        }\n    }\n}\n\nif ($found) {\n    # Replace \"PF\" with \"Cuest.\" in place; formatting (bold + accent2\n    # orange color) carries over since we're reusing the same range.\n    $rng.Text = \"Cuest.\"\n\n    # Move the \"_GoBack\" bookmark to right after the text we just typed.\n    if ($d.Bookmarks.Exists(\"_GoBack\")) {\n        $d.Bookmarks.Item(\"_GoBack\").Delete()\n    }\n    $bmRange = $d.Range($rng.End, $rng.End)\n    $d.Bookmarks.Add(\"_GoBack\", $bmRange)\n}\n"}
